# Rehearsal application documentation update
# Adds a "Possible future improvements:" section with a 3-item numbered
# list, inserted among the trailing blank paragraphs before the sectPr.

$d = $word.ActiveDocument

# Locate the run of trailing empty paragraphs that precede the final
# section break. We insert our new content before the 3rd-from-last one,
# so 3 blank paragraphs remain after the new block (mirroring the source
# edit, which leaves the new content flanked by blank paragraphs).
$total = $d.Paragraphs.Count
$insertBeforeIndex = $total - 2

$targetPara = $d.Paragraphs.Item($insertBeforeIndex)
$insPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)

# Insert the heading paragraph plus the three list-item paragraphs in one
# shot (using `\r` as the paragraph mark), then fix up styling/numbering.
$block = "Possible future improvements:`r" + `
         "Create different types of users (admin, regular). Admins are able to give to other regular users admin rights.`r" + `
         "Implement notifications for important notes.`r" + `
         "Implement a chat.`r"
$insPoint.InsertBefore($block)

# Indices (1-based) of the 4 newly inserted paragraphs.
$headingIndex = $insertBeforeIndex
$item1Index = $insertBeforeIndex + 1
$item2Index = $insertBeforeIndex + 2
$item3Index = $insertBeforeIndex + 3

# Turn the three list paragraphs into a single freshly-numbered list
# (decimal "1.", "2.", "3."), all sharing one list/numbering definition.
# Apply the "List Paragraph" style to all three *before* touching the
# numbering, otherwise a later style assignment clobbers the numPr that
# was just applied.
$item1 = $d.Paragraphs.Item($item1Index)
$item1.Style = "List Paragraph"
$item2 = $d.Paragraphs.Item($item2Index)
$item2.Style = "List Paragraph"
$item3 = $d.Paragraphs.Item($item3Index)
$item3.Style = "List Paragraph"

$item1.Range.ListFormat.ApplyNumberDefault()

$rest = $d.Range($d.Paragraphs.Item($item2Index).Range.Start, $d.Paragraphs.Item($item3Index).Range.End)
$rest.ListFormat.ApplyNumberDefault()

Write-Host "Inserted 'Possible future improvements' section with 3 list items."
